$d = $word.ActiveDocument

# --- Fix 1: merge the "Markt Survey C:" run split into a single run ---
# (Original had " Survey " / "C" / ":" as three separate runs; final has one
#  run with " Survey C:". Replacing the literal text with itself collapses
#  the runs into one.)
$d.Content.Find.Execute(" Survey C:", $true, $false, $false, $false, $false, $true, 1, $false, " Survey C:", 2) | Out-Null

# --- Fix 2: append the new "Cross reference sweep 2/3", "Check markets on
#     Tanzania Map" and "Additional cleaning steps" sections after the last
#     paragraph ("Repeat for markets") ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$insertionPoint = $d.Paragraphs.Last

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p/><w:p/><w:p>
      <w:r>
        <w:t>Cross reference sweep 2:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Populate missing names in </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>eithr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> the market or village column based on a completed market or village name (</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>i.e.</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> backfilling)</w:t>
      </w:r>
    </w:p><w:p/><w:p>
      <w:r>
        <w:t xml:space="preserve">Cross reference sweep </w:t>
      </w:r>
      <w:r>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Populate remaining missing entries based on original entry. </w:t>
      </w:r>
    </w:p><w:p/><w:p/><w:p>
      <w:r>
        <w:t xml:space="preserve">Check markets on Tanzania Map: </w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Use names cleaning script and plot locations of corrected market </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>names</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Check that markets under same names are in the same location – corrections listed on Tanzania data presentation. </w:t>
      </w:r>
    </w:p><w:p/><w:p/><w:p>
      <w:r>
        <w:t>Additional cleaning steps:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Check ward </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>names</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Clean </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>gps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> coordinates – some have </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>mistakes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p><w:p/>
'@

$insertionPoint.Range.InsertXML($newXml)

Write-Output "paragraphs now: $($d.Paragraphs.Count)"
